$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new cells on the existing last row (row 9)
$ws.Range("X9").Value = -1.7200020000000222
$ws.Range("Y9").Value = "Down"

# Add a brand new row (row 10) with a full set of values
$ws.Range("A10").Value = 42653.880173611113
$ws.Range("A10").NumberFormat = "m/d/yy h:mm"

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 38
$ws.Range("E10").Value = 8988
$ws.Range("F10").Value = 1940
$ws.Range("G10").Value = 70
$ws.Range("H10").Value = 28
$ws.Range("I10").Value = 83
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 13387
$ws.Range("L10").Value = 196
$ws.Range("M10").Value = 79
$ws.Range("N10").Value = 69
$ws.Range("O10").Value = 13
$ws.Range("P10").Value = "Noun"
$ws.Range("Q10").Value = 44.409433632991338
$ws.Range("R10").Value = 1.8

$ws.Range("S10").Value = 0.0926
$ws.Range("S10").NumberFormat = "0.00%"

$ws.Range("T10").Value = -0.0094
$ws.Range("T10").NumberFormat = "0.00%"

$ws.Range("U10").Value = 5.87
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = 0

# The new row shifts the "best fit" auto width of each populated column
# very slightly (this mirrors Excel's own auto bestFit re-measurement
# after the data changed).
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.6666666666666667
$ws.Columns.Item(3).ColumnWidth = 8.3333333333333333
$ws.Columns.Item(4).ColumnWidth = 11.3333333333333333
$ws.Columns.Item(5).ColumnWidth = 8.6666666666666667
$ws.Columns.Item(6).ColumnWidth = 11.3333333333333333
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(10).ColumnWidth = 19.8333333333333333
$ws.Columns.Item(11).ColumnWidth = 9.5
$ws.Columns.Item(12).ColumnWidth = 13.5
$ws.Columns.Item(13).ColumnWidth = 13.8333333333333333

$wb.Save()
